$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.127.02"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.108.60"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  -0.96%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "350.71"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +4.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5172"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -1.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4504"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -1.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.86"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -4.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08974"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.176"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.73"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +5.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.108.73"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.768"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.151"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.63"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +2.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001151"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -2.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.003"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -0.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.61"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +6.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06676"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("E21").Value = "  -0.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.244"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.226.06"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -1.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.90"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.345"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.356.62"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.15"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.565"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +1.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.96"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.49"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.189"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -2.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1069"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.655"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +1.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.278"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -1.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.964"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.960"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +1.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.22"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -2.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02592"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -0.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06857"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2318"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -0.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.56"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -0.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6842"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.37"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +2.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6436"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.298"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.677"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("E48").Value = "  +4.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "83.80"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +0.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07230"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +0.62%  "
